$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 999 (old rows 999-1070 shift down to 1002-1073)
$ws.Rows("999:1001").Insert()

# Fill in the 3 new rows with fresh data (same shape as the rest of the table)

# Row 999
$ws.Range("A999").Value = 4
$ws.Range("B999").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C999").Value = "Los Lagos"
$ws.Range("D999").Value = 45265
$ws.Range("E999").Value = 10
$ws.Range("F999").Value = 100112004
$ws.Range("G999").Value = "Cebolla"
$ws.Range("H999").Value = "Morada(o)"
$ws.Range("I999").Value = "Primera"
$ws.Range("J999").Value = 300
$ws.Range("K999").Value = 18000
$ws.Range("L999").Value = 19000
$ws.Range("M999").Value = 18500
$ws.Range("N999").Value = "`$/malla 18 kilos"
$ws.Range("O999").Value = "Perú"
$ws.Range("P999").Value = 1028
$ws.Range("Q999").Value = 18
$ws.Range("R999").Value = "Hortaliza"

# Row 1000
$ws.Range("A1000").Value = 4
$ws.Range("B1000").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C1000").Value = "Los Lagos"
$ws.Range("D1000").Value = 45265
$ws.Range("E1000").Value = 10
$ws.Range("F1000").Value = 100112004
$ws.Range("G1000").Value = "Cebolla"
$ws.Range("H1000").Value = "Sin especificar"
$ws.Range("I1000").Value = "1a (cosecha)"
$ws.Range("J1000").Value = 750
$ws.Range("K1000").Value = 19000
$ws.Range("L1000").Value = 19000
$ws.Range("M1000").Value = 19000
$ws.Range("N1000").Value = "`$/malla 17 kilos"
$ws.Range("O1000").Value = "Región de O'Higgins"
$ws.Range("P1000").Value = 1118
$ws.Range("Q1000").Value = 17
$ws.Range("R1000").Value = "Hortaliza"

# Row 1001
$ws.Range("A1001").Value = 4
$ws.Range("B1001").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C1001").Value = "Los Lagos"
$ws.Range("D1001").Value = 45265
$ws.Range("E1001").Value = 10
$ws.Range("F1001").Value = 100112004
$ws.Range("G1001").Value = "Cebolla"
$ws.Range("H1001").Value = "Sin especificar"
$ws.Range("I1001").Value = "Primera"
$ws.Range("J1001").Value = 750
$ws.Range("K1001").Value = 19000
$ws.Range("L1001").Value = 19000
$ws.Range("M1001").Value = 19000
$ws.Range("N1001").Value = "`$/malla 18 kilos"
$ws.Range("O1001").Value = "Perú"
$ws.Range("P1001").Value = 1056
$ws.Range("Q1001").Value = 18
$ws.Range("R1001").Value = "Hortaliza"

# Ensure date style is preserved on D column for the new rows (matching the rest of the table)
$ws.Range("D999:D1001").NumberFormat = $ws.Range("D1002").NumberFormat
